$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.590.74"
$ws.Range("E2").Value = "  -0.18%  "
Set-TextValue "D3" "1.593.70"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "210.71"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue "D9" "0.244"
$ws.Range("E9").Value = "  -0.97%  "
Set-TextValue "D10" "19.33"
$ws.Range("E10").Value = "  -1.59%  "
Set-TextValue "D11" "0.0835"
$ws.Range("E11").Value = "  +0.13%  "
Set-TextValue "D12" "1.817.67"
$ws.Range("E12").Value = "  +0.30%  "
Set-TextValue "D13" "1.584.34"
$ws.Range("E13").Value = "  -0.48%  "
Set-TextValue "D15" "0.519"
$ws.Range("E15").Value = "  -0.62%  "
Set-TextValue "D16" "64.31"
$ws.Range("E16").Value = "  -0.80%  "
Set-TextValue "D17" "26.576.35"
Set-TextValue "D18" "0.0₃0728"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.05%  "
Set-TextValue "D20" "207.53"
$ws.Range("E20").Value = "  -0.47%  "
Set-TextValue "D21" "6.92"
$ws.Range("E21").Value = "  +2.80%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  -0.34%  "
Set-TextValue "D25" "144.99"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("E26").Value = "  +0.01%  "
Set-TextValue "D27" "7.08"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").Value = "  -0.60%  "
Set-TextValue "D30" "0.0503"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("E32").Value = "  -0.09%  "
Set-TextValue "D33" "0.652"
$ws.Range("E33").Value = "  -1.53%  "
Set-TextValue "D34" "2.91"
$ws.Range("E34").Value = "  +0.56%  "
Set-TextValue "D35" "1.276.78"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  +0.09%  "
Set-TextValue "D41" "5.45"
$ws.Range("E41").Value = "  +1.69%  "
Set-TextValue "D42" "2.20"
$ws.Range("E42").Value = "  +1.46%  "
Set-TextValue "D43" "0.784"
$ws.Range("E43").Value = "  -1.03%  "
Set-TextValue "D44" "63.79"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("E45").Value = "  +9.94%  "
Set-TextValue "D46" "1.729.58"
$ws.Range("E46").Value = "  +0.22%  "
Set-TextValue "D47" "89.29"
$ws.Range("E47").Value = "  -0.45%  "
Set-TextValue "D48" "1.58"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  +3.55%  "
Set-TextValue "D51" "0.0505"
$ws.Range("E51").Value = "  +0.35%  "
